# Auto-generated Excel COM-interop script to apply the Bahamut_Profits market-data refresh diff.
# For each affected sheet/row, update columns H-N to the new scraped values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 103.21875
$ws.Range("J33").Value = 105.75
$ws.Range("L33").Value = 105.75
$ws.Range("N33").Value = -563.75
$ws.Range("H74").Value = 3248219.5
$ws.Range("I74").Value = 3464534
$ws.Range("J74").Value = 3500
$ws.Range("K74").Value = 3464534
$ws.Range("L74").Value = 3500
$ws.Range("M74").Value = -3463598
$ws.Range("N74").Value = -5372
$ws.Range("H77").Value = 3248219.5
$ws.Range("I77").Value = 3464534
$ws.Range("J77").Value = 3500
$ws.Range("K77").Value = 17322670
$ws.Range("L77").Value = 17500
$ws.Range("M77").Value = -17317990
$ws.Range("N77").Value = -26860
$ws.Range("H86").Value = 2556.121
$ws.Range("I86").Value = 2200.9375
$ws.Range("J86").Value = 2890.4119
$ws.Range("K86").Value = 2200.9375
$ws.Range("L86").Value = 2890.4119
$ws.Range("M86").Value = -1077.9375
$ws.Range("N86").Value = -5136.4119
$ws.Range("H89").Value = 2556.121
$ws.Range("I89").Value = 2200.9375
$ws.Range("J89").Value = 2890.4119
$ws.Range("K89").Value = 11004.6875
$ws.Range("L89").Value = 14452.0595
$ws.Range("M89").Value = -5388.6875
$ws.Range("N89").Value = -25684.0595
$ws.Range("H100").Value = 1610.5
$ws.Range("I100").Value = 1554.8572
$ws.Range("J100").Value = 2000
$ws.Range("K100").Value = 1554.8572
$ws.Range("L100").Value = 2000
$ws.Range("M100").Value = -1013.8572
$ws.Range("N100").Value = -3082
$ws.Range("H129").Value = 1482553.5
$ws.Range("I129").Value = 285.1
$ws.Range("J129").Value = 2470732.5
$ws.Range("K129").Value = 855.3000000000001
$ws.Range("L129").Value = 7412197.5
$ws.Range("M129").Value = 4144.7
$ws.Range("N129").Value = -7422197.5
$ws.Range("H132").Value = 1207.7606
$ws.Range("I132").Value = 1274.9193
$ws.Range("J132").Value = 745.1111
$ws.Range("K132").Value = 3824.7579
$ws.Range("L132").Value = 2235.3333
$ws.Range("M132").Value = -1294.7579
$ws.Range("N132").Value = -7295.3333
$ws.Range("H135").Value = 697.4394
$ws.Range("I135").Value = 384.5263
$ws.Range("K135").Value = 3460.7367
$ws.Range("M135").Value = -925.7366999999999
$ws.Range("H137").Value = 845
$ws.Range("I137").Value = 745.63635
$ws.Range("J137").Value = 894.6818
$ws.Range("K137").Value = 2236.90905
$ws.Range("L137").Value = 2684.0454
$ws.Range("M137").Value = 313.0909499999998
$ws.Range("N137").Value = -7784.0454
$ws.Range("H138").Value = 1541.4
$ws.Range("I138").Value = 755
$ws.Range("J138").Value = 2359.898
$ws.Range("K138").Value = 2265
$ws.Range("L138").Value = 7079.694
$ws.Range("M138").Value = 2875
$ws.Range("N138").Value = -17359.694
$ws.Range("H141").Value = 2422.2559
$ws.Range("I141").Value = 987.0571
$ws.Range("K141").Value = 2961.1713
$ws.Range("M141").Value = 2218.8287

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4092.43
$ws.Range("I32").Value = 3687.6
$ws.Range("J32").Value = 6386.467
$ws.Range("K32").Value = 3687.6
$ws.Range("L32").Value = 6386.467
$ws.Range("M32").Value = -3400.6
$ws.Range("N32").Value = -6960.467
$ws.Range("H55").Value = 251361.22
$ws.Range("J55").Value = 251361.22
$ws.Range("L55").Value = 251361.22
$ws.Range("N55").Value = -251991.22
$ws.Range("H61").Value = 860.4524
$ws.Range("I61").Value = 720.4167
$ws.Range("J61").Value = 1700.6666
$ws.Range("K61").Value = 720.4167
$ws.Range("L61").Value = 1700.6666
$ws.Range("M61").Value = -508.4167
$ws.Range("N61").Value = -2124.6666
$ws.Range("H63").Value = 1996.0714
$ws.Range("I63").Value = 2007.36
$ws.Range("J63").Value = 1902
$ws.Range("K63").Value = 2007.36
$ws.Range("L63").Value = 1902
$ws.Range("M63").Value = -1321.36
$ws.Range("N63").Value = -3274
$ws.Range("H66").Value = 1996.0714
$ws.Range("I66").Value = 2007.36
$ws.Range("J66").Value = 1902
$ws.Range("K66").Value = 10036.8
$ws.Range("L66").Value = 9510
$ws.Range("M66").Value = -6604.799999999999
$ws.Range("N66").Value = -16374
$ws.Range("H74").Value = 899.2641599999999
$ws.Range("I74").Value = 895.57776
$ws.Range("J74").Value = 920
$ws.Range("K74").Value = 895.57776
$ws.Range("L74").Value = 920
$ws.Range("M74").Value = -21.57776000000001
$ws.Range("N74").Value = -2668
$ws.Range("H77").Value = 899.2641599999999
$ws.Range("I77").Value = 895.57776
$ws.Range("J77").Value = 920
$ws.Range("K77").Value = 4477.8888
$ws.Range("L77").Value = 4600
$ws.Range("M77").Value = -109.8887999999997
$ws.Range("N77").Value = -13336
$ws.Range("H132").Value = 1008.92
$ws.Range("I132").Value = 906.7895
$ws.Range("J132").Value = 1332.3334
$ws.Range("K132").Value = 2720.3685
$ws.Range("L132").Value = 3997.0002
$ws.Range("M132").Value = -190.3685
$ws.Range("N132").Value = -9057.0002
$ws.Range("H136").Value = 860.4524
$ws.Range("I136").Value = 720.4167
$ws.Range("J136").Value = 1700.6666
$ws.Range("K136").Value = 2161.2501
$ws.Range("L136").Value = 5101.9998
$ws.Range("M136").Value = 388.7498999999998
$ws.Range("N136").Value = -10201.9998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2038.2858
$ws.Range("I86").Value = 1904.081
$ws.Range("K86").Value = 1904.081
$ws.Range("M86").Value = -781.0809999999999
$ws.Range("H89").Value = 2038.2858
$ws.Range("I89").Value = 1904.081
$ws.Range("K89").Value = 9520.404999999999
$ws.Range("M89").Value = -3904.404999999999
$ws.Range("H134").Value = 18144.623
$ws.Range("I134").Value = 1520.8269
$ws.Range("J134").Value = 114193.22
$ws.Range("K134").Value = 4562.4807
$ws.Range("L134").Value = 342579.66
$ws.Range("M134").Value = -2027.4807
$ws.Range("N134").Value = -347649.66

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2665.2307
$ws.Range("I31").Value = 2513.875
$ws.Range("J31").Value = 3357.1428
$ws.Range("K31").Value = 2513.875
$ws.Range("L31").Value = 3357.1428
$ws.Range("M31").Value = -2218.875
$ws.Range("N31").Value = -3947.1428
$ws.Range("H34").Value = 2665.2307
$ws.Range("I34").Value = 2513.875
$ws.Range("J34").Value = 3357.1428
$ws.Range("K34").Value = 2513.875
$ws.Range("L34").Value = 3357.1428
$ws.Range("M34").Value = -2311.875
$ws.Range("N34").Value = -3761.1428
$ws.Range("H58").Value = 2471.0352
$ws.Range("I58").Value = 610.38464
$ws.Range("J58").Value = 6502.4443
$ws.Range("K58").Value = 610.38464
$ws.Range("L58").Value = 6502.4443
$ws.Range("M58").Value = -407.38464
$ws.Range("N58").Value = -6908.4443
$ws.Range("H132").Value = 1372.8369
$ws.Range("I132").Value = 858.4035
$ws.Range("J132").Value = 2210.6287
$ws.Range("K132").Value = 2575.2105
$ws.Range("L132").Value = 6631.886100000001
$ws.Range("M132").Value = -45.21050000000014
$ws.Range("N132").Value = -11691.8861
$ws.Range("H134").Value = 1240
$ws.Range("I134").Value = 1177.3489
$ws.Range("K134").Value = 3532.0467
$ws.Range("M134").Value = -997.0466999999999
$ws.Range("H136").Value = 2471.0352
$ws.Range("I136").Value = 610.38464
$ws.Range("J136").Value = 6502.4443
$ws.Range("K136").Value = 1831.15392
$ws.Range("L136").Value = 19507.3329
$ws.Range("M136").Value = 718.84608
$ws.Range("N136").Value = -24607.3329

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H45").Value = 978.3
$ws.Range("J45").Value = 1090.4286
$ws.Range("L45").Value = 3271.2858
$ws.Range("N45").Value = -4335.2858
$ws.Range("H92").Value = 501.5
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 501.5
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 1504.5
$ws.Range("M92").ClearContents()
$ws.Range("N92").Value = -4000.5
$ws.Range("H131").Value = 31298.242
$ws.Range("I131").Value = 112472.336
$ws.Range("J131").Value = 18481.281
$ws.Range("K131").Value = 337417.008
$ws.Range("L131").Value = 55443.84299999999
$ws.Range("M131").Value = -332377.008
$ws.Range("N131").Value = -65523.84299999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1561.1842
$ws.Range("I132").Value = 1388.5319
$ws.Range("J132").Value = 1841
$ws.Range("K132").Value = 4165.5957
$ws.Range("L132").Value = 5523
$ws.Range("M132").Value = -1635.5957
$ws.Range("N132").Value = -10583

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1082.4783
$ws.Range("I22").Value = 380.7857
$ws.Range("J22").Value = 1389.4688
$ws.Range("K22").Value = 380.7857
$ws.Range("L22").Value = 1389.4688
$ws.Range("M22").Value = -85.78570000000002
$ws.Range("N22").Value = -1979.4688
$ws.Range("H27").Value = 1082.4783
$ws.Range("I27").Value = 380.7857
$ws.Range("J27").Value = 1389.4688
$ws.Range("K27").Value = 380.7857
$ws.Range("L27").Value = 1389.4688
$ws.Range("M27").Value = -273.7857
$ws.Range("N27").Value = -1603.4688
$ws.Range("H55").Value = 247.2963
$ws.Range("I55").Value = 169.6923
$ws.Range("J55").Value = 319.35715
$ws.Range("K55").Value = 169.6923
$ws.Range("L55").Value = 319.35715
$ws.Range("M55").Value = 3.307700000000011
$ws.Range("N55").Value = -665.35715
$ws.Range("H132").Value = 2092.6345
$ws.Range("I132").Value = 1939.8158
$ws.Range("J132").Value = 2507.4285
$ws.Range("K132").Value = 5819.4474
$ws.Range("L132").Value = 7522.2855
$ws.Range("M132").Value = -3289.4474
$ws.Range("N132").Value = -12582.2855
$ws.Range("H136").Value = 1374.5
$ws.Range("I136").Value = 778.597
$ws.Range("J136").Value = 5004.091
$ws.Range("K136").Value = 2335.791
$ws.Range("L136").Value = 15012.273
$ws.Range("M136").Value = 214.2089999999998
$ws.Range("N136").Value = -20112.273

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 858.54346
$ws.Range("I132").Value = 688.2059
$ws.Range("J132").Value = 1341.1666
$ws.Range("K132").Value = 2064.6177
$ws.Range("L132").Value = 4023.4998
$ws.Range("M132").Value = 465.3822999999998
$ws.Range("N132").Value = -9083.4998
$ws.Range("H136").Value = 832.2683
$ws.Range("I136").Value = 897.5484
$ws.Range("J136").Value = 629.9
$ws.Range("K136").Value = 2692.6452
$ws.Range("L136").Value = 1889.7
$ws.Range("M136").Value = -142.6451999999999
$ws.Range("N136").Value = -6989.7
